# BIIBBagTrade20.xlsx - fix the 20 minute trade row (row 3):
#   - recompute the trade using the actual Buy/Sell prices recorded
#     for this run instead of the placeholder values
#   - Principle (A3) rolls forward from the new Price Change % (F3)
#   - Date (G3) reflects when the trade actually closed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkSheet 1")

$ws.Range("A3").Value = 9983
$ws.Range("C3").Value = 309.02999999999997
$ws.Range("D3").Value = 309.55
$ws.Range("F3").Value = 0.17
$ws.Range("G3").Value = 42608.637824074074
